$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149, shifting existing rows 149:170 down to 150:171.
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,N,Q,R repeat the constant values used throughout
# this "Ciboulette" / "Vega Modelo de Temuco" block.
$ws.Cells.Item(149, 1).Value = 10
$ws.Cells.Item(149, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(149, 3).Value = "La Araucanía"
$ws.Cells.Item(149, 4).Value = 44505
$ws.Cells.Item(149, 4).Style = $ws.Cells.Item(150, 4).Style
$ws.Cells.Item(149, 4).NumberFormat = $ws.Cells.Item(150, 4).NumberFormat
$ws.Cells.Item(149, 5).Value = 9
$ws.Cells.Item(149, 6).Value = 100112039
$ws.Cells.Item(149, 7).Value = "Ciboulette"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 65
$ws.Cells.Item(149, 11).Value = 6000
$ws.Cells.Item(149, 12).Value = 6000
$ws.Cells.Item(149, 13).Value = 6000
$ws.Cells.Item(149, 14).Value = "$/docena de atados"
$ws.Cells.Item(149, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(149, 16).Value = 2000
$ws.Cells.Item(149, 17).Value = 3
$ws.Cells.Item(149, 18).Value = "Hortaliza"
